$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 1121, pushing existing data (old 1121..1201) down to (1123..1203)
$ws.Rows.Item(1121).Resize(2).Insert()

# New row 1121: based on old row 1121 (now at 1123), but with updated fields
$ws.Range("A1121:T1121").Value2 = $ws.Range("A1123:T1123").Value2
$ws.Range("D1121").Value2 = 44783
$ws.Range("M1121").Value2 = 150
$ws.Range("N1121").Value2 = 28000
$ws.Range("O1121").Value2 = 28000
$ws.Range("P1121").Value2 = 28000
$ws.Range("R1121").Value2 = "Provincia de Melipilla"
$ws.Range("S1121").Value2 = 4000

# New row 1122: based on old row 1122 (now at 1124), but with updated fields
$ws.Range("A1122:T1122").Value2 = $ws.Range("A1124:T1124").Value2
$ws.Range("D1122").Value2 = 44783
$ws.Range("L1122").Value2 = "Segunda"
$ws.Range("M1122").Value2 = 100
$ws.Range("N1122").Value2 = 18000
$ws.Range("O1122").Value2 = 18000
$ws.Range("P1122").Value2 = 18000
$ws.Range("R1122").Value2 = "Provincia de Melipilla"
$ws.Range("S1122").Value2 = 2571
